$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset for rows 2..31 (timestamp, label, ax, ay, az, gx, gy, gz).
# 9 new sensor-reading rows were prepended right after the header (pushing the
# previous 20 data rows down by 9), and 1 new row was appended at the end, so
# the sheet grows from 20 data rows (A1:H21) to 30 data rows (A1:H31). The
# timestamp/label columns keep following the original 0,100,200.. / "falling"
# pattern for every row.
$data = @(
    @(0, "falling", -2.203232884407044, 3.565518021583557, 2.421210885047912, -0.1018617823719978, 0.0229074470698833, 0.0226020142436027),
    @(100, "falling", -2.176557183265686, 3.694440901279449, 2.333479553461074, -0.0087048299610614, -0.0256563406437635, -0.0074830991216003),
    @(200, "falling", -2.326471328735352, 3.527873754501343, 2.532370328903198, -0.0534507073462009, -0.009010262787342, 0.001527163083665),
    @(300, "falling", -2.219938099384308, 3.655792444944382, 2.356739014387131, -0.0387899428606033, 0.0113010071218013, 0.0397062413394451),
    @(400, "falling", -2.15459930896759, 3.624111294746399, 2.422300338745117, -0.011148290708661, -0.08338310569524759, 0.0064140851609408),
    @(500, "falling", -2.086719453334808, 3.556297302246094, 2.323985010385513, -0.0229074470698833, -0.0862847194075584, 0.0215329993516206),
    @(600, "falling", -2.057324409484864, 3.517749786376953, 2.224772691726685, -0.0293215326964855, -0.0226020142436027, -0.0091629782691597),
    @(700, "falling", -2.064002573490143, 3.512480080127716, 2.171477824449539, -0.0369573459029197, -0.0378736443817615, -0.0192422550171613),
    @(800, "falling", -2.053778767585754, 3.516274869441986, 2.167593479156494, 0.0308486949652433, -0.0488692186772823, -0.0198531206697225),
    @(900, "falling", -2.025566756725311, 3.52062651515007, 2.27691987156868, -0.0397062413394451, -0.0024434609804302, 0.0332921557128429),
    @(1000, "falling", -2.230706214904786, 3.561713695526123, 2.031704187393189, -0.0200058370828628, -0.0035124751739203, 0.0421497002243995),
    @(1100, "falling", -2.475497364997865, 3.523229420185089, 2.537566900253297, -0.0039706239476799, 0.0705549344420433, 0.0328340083360672),
    @(1200, "falling", -2.48767375946045, 3.061535835266112, 3.86157149076462, -0.0146607663482427, 0.2353358417749405, 0.07696902006864539),
    @(1300, "falling", -2.442349374294281, 2.902286112308502, 4.137762367725372, 0.0704022198915481, 0.2420553565025329, 0.08109235763549801),
    @(1400, "falling", -2.809046030044556, 1.600962877273557, 4.795935153961183, 0.6050620079040527, 0.2802344262599945, 0.07590000331401819),
    @(1500, "falling", -3.426159977912904, -0.2882512211799662, 5.237384021282194, 0.1878410577774047, 0.7900014519691467, 0.304669052362442),
    @(1600, "falling", -4.615099787712106, 0.4410536289215286, 4.911977410316469, 0.3132211565971374, -0.3736968040466308, 0.0494800843298435),
    @(1700, "falling", -5.570275843143422, 5.792195498943334, 7.507518291473446, 0.6346889734268188, -0.973566472530365, 0.2535090744495392),
    @(1800, "falling", 4.738878250122034, 7.207733154296871, 21.9623451232909, -2.739883422851562, 3.34677791595459, -2.103209018707275),
    @(1900, "falling", -3.471131086349493, 6.420480489730835, -4.152171969413772, -3.001944541931152, 0.3590360581874847, 1.426217675209045),
    @(2000, "falling", -2.357546925544736, 3.453876137733455, 3.317040443420424, -0.064446285367012, -0.0032070425804704, -1.773189067840576),
    @(2100, "falling", -0.1309916377067535, 3.795689940452578, 3.876870155334465, -0.5893322229385376, -0.0039706239476799, -0.3778201639652252),
    @(2200, "falling", -0.2815589904785198, 4.85230040550232, 2.205311059951784, 0.3320052623748779, -0.7985535860061646, 0.4050036668777466),
    @(2300, "falling", -1.535586237907411, 4.95532476902008, 1.937351673841474, 0.0734565481543541, -1.579086661338806, -1.274570345878601),
    @(2400, "falling", -1.377771139144897, 4.249351501464844, 1.058028712868691, -0.2518292069435119, -0.9622654914855956, -0.384845107793808),
    @(2500, "falling", 0.1868795156478872, 3.352623224258432, 1.056536458432678, -0.0522289797663688, -0.2220495194196701, -0.2014328092336654),
    @(2600, "falling", -0.04321670532226736, 5.682518005371086, 1.995282649993894, 0.0580321997404098, -0.2347249686717987, 0.4702135324478149),
    @(2700, "falling", -0.4393689632415791, 3.92675977945327, 1.474137753248212, 0.08491026610136029, 0.1505782902240753, -0.0226020142436027),
    @(2800, "falling", -0.5461759567260746, 4.010827064514163, 1.266485691070556, 0.0704022198915481, 0.0319177098572254, 0.1357648074626922),
    @(2900, "falling", -0.5815373659133911, 4.265253961086274, 1.277935206890107, -0.0236710291355848, 0.039248090237379, 0.0468839071691036)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
